# "Generate Report for Handback"
# Updates the localization-status workbook to reflect that the zh-cn and
# de-de handback packages have come back from translation: the Status
# column flips from "Ready for handoff" to "Handed back: in sync with
# en-US", the Latest Target File / Latest Handback File / Latest Handback
# DateTime columns get populated (with a hyperlink on the target file),
# and the relevant columns are widened so the new text is readable.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1c5eae0cdc6f89a1b511c27272241337a4fde9d6/e2e/879e77d3-ab9b-4ed8-8ae1-6a51bd3c5903.md"
$mdDisplay = "879e77d3-ab9b-4ed8-8ae1-6a51bd3c5903.md"
$newStatus = "Handed back: in sync with en-US"

# The "Status" text is shared across the Overview rollup columns (zh-cn /
# de-de) and each language sheet's own Status cell - update every place it
# shows up.
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# --- zh-cn sheet -----------------------------------------------------
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("J2").Value = "879e77d3-ab9b-4ed8-8ae1-6a51bd3c5903.8044a48629369dd80ffd6be8418e73f80b4fdbe6.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-29 23:02:30"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl, "", "", $mdDisplay)

# --- de-de sheet -------------------------------------------------------
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("J2").Value = "879e77d3-ab9b-4ed8-8ae1-6a51bd3c5903.8044a48629369dd80ffd6be8418e73f80b4fdbe6.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-29 23:02:38"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl, "", "", $mdDisplay)

# --- widen columns to fit the new, longer content ----------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664
